$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 2023-10-25 (45224) to 2023-11-03 (45233)
$newDate = (Get-Date -Year 2023 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
$ws.Range("C6").Value = $newDate
